$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("C1,J1,M1,O1,P1,Q1")
foreach ($area in $r.Areas) {
    $area.Interior.Color = 65535
}
